$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 68
$ws.Range("I2").Value = 167
$ws.Range("J2").Value = 689
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 181
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 140
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 88
$ws.Range("T2").Value = 135
$ws.Range("U2").Value = 11
$ws.Range("V2").Value = 1084
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1110
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 20
$ws.Range("AA2").Value = 9
